$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = 4.621579
$ws.Range("N2").Value = 13.864737
$ws.Range("O2").Value = 0.1778708528171788
$ws.Range("P2").Value = 0.1778708528171788
$ws.Range("Q2").Value = 0.3153981183286666
$ws.Range("R2").Value = 2.838583064958
$ws.Range("S2").Value = 0.1766833786845485
$ws.Range("T2").Value = 0.1766833786845485
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5971062807549863
$ws.Range("P3").Value = 0.5971062807549863
$ws.Range("R3").Value = 9.529024849694
$ws.Range("S3").Value = 0.5931199712973236
$ws.Range("T3").Value = 0.5931199712973236
$ws.Range("O4").Value = 0.2250228664278349
$ws.Range("P4").Value = 0.2250228664278349
$ws.Range("R4").Value = 3.591066707972
$ws.Range("S4").Value = 0.2235206032469864
$ws.Range("T4").Value = 0.2235206032469864
$ws.Range("J5").Value = 0.006676046771141624
$ws.Range("M5").Value = 4.621579
$ws.Range("N5").Value = 13.864737
$ws.Range("O5").Value = 0.1778708528171788
$ws.Range("P5").Value = 0.1778708528171788
$ws.Range("Q5").Value = 0.002119764234666666
$ws.Range("R5").Value = 0.019077878112
$ws.Range("S5").Value = 0.001187474132630334
$ws.Range("T5").Value = 0.001187474132630334
$ws.Range("J6").Value = 0.006676046771141624
$ws.Range("N6").Value = 46.543441
$ws.Range("O6").Value = 0.5971062807549863
$ws.Range("P6").Value = 0.5971062807549863
$ws.Range("Q6").Value = 0.007115974979555556
$ws.Range("S6").Value = 0.00398630945766271
$ws.Range("T6").Value = 0.00398630945766271
$ws.Range("J7").Value = 0.006676046771141624
$ws.Range("O7").Value = 0.2250228664278349
$ws.Range("P7").Value = 0.2250228664278349
$ws.Range("S7").Value = 0.00150226318084858
$ws.Range("T7").Value = 0.00150226318084858
